$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 9.271299999999998
$ws.Range("D4").Value = -7.7293
$ws.Range("D7").Value = -7.690800000000005
$ws.Range("A9").Value = -22.10599999999999
$ws.Range("B9").Value = 6.666300000000007
$ws.Range("C9").Value = -12.09410000000001
$ws.Range("B11").Value = 5.2858
$ws.Range("D11").Value = -7.683799999999992
$ws.Range("A13").Value = -22.43950000000001
$ws.Range("D15").Value = -8.309699999999998
$ws.Range("A16").Value = -21.7353
$ws.Range("B16").Value = 4.213299999999995
$ws.Range("A18").Value = -22.26320000000001
$ws.Range("A20").Value = -20.13149999999999
$ws.Range("C22").Value = -12.8378
$ws.Range("B23").Value = 9.407599999999995
$ws.Range("D23").Value = -7.054399999999999
$ws.Range("B24").Value = 5.2687
$ws.Range("A26").Value = -21.4186
$ws.Range("B26").Value = 4.781600000000001
$ws.Range("A27").Value = -22.022
$ws.Range("C27").Value = -12.88849999999999
$ws.Range("A29").Value = -21.65009999999997
$ws.Range("C29").Value = -12.23570000000001
$ws.Range("D30").Value = -7.217099999999999
$ws.Range("C32").Value = -12.6647
$ws.Range("B34").Value = 9.524700000000003
$ws.Range("A35").Value = -21.717
$ws.Range("B35").Value = 3.813099999999997
$ws.Range("A36").Value = -20.69339999999999
$ws.Range("C37").Value = -14.38229999999998
$ws.Range("C38").Value = -11.58870000000001
$ws.Range("D38").Value = -7.414199999999999
$ws.Range("C39").Value = -12.70600000000002
$ws.Range("D39").Value = -7.652400000000004
$ws.Range("C41").Value = -12.4254
$ws.Range("D43").Value = -7.549400000000009
$ws.Range("B44").Value = 4.831200000000004
$ws.Range("A45").Value = -21.68409999999999
$ws.Range("C45").Value = -13.61469999999998
$ws.Range("D47").Value = -7.378300000000003
$ws.Range("B48").Value = 6.825800000000003
$ws.Range("C48").Value = -12.51400000000001
$ws.Range("B49").Value = 4.924
$ws.Range("C51").Value = -11.21129999999999
$ws.Range("B52").Value = 5.4836
$ws.Range("A55").Value = -22.06229999999999
$ws.Range("C56").Value = -12.71849999999999
$ws.Range("A57").Value = -22.0938
$ws.Range("C57").Value = -12.58649999999999
$ws.Range("C61").Value = -14.35909999999998
$ws.Range("C64").Value = -10.44
$ws.Range("B66").Value = 4.879799999999995
$ws.Range("B67").Value = 4.868099999999998
$ws.Range("A69").Value = -21.55819999999998
$ws.Range("B73").Value = 9.220100000000008
$ws.Range("C75").Value = -12.04139999999998
$ws.Range("D75").Value = -7.4918
$ws.Range("A76").Value = -19.96309999999999
$ws.Range("A78").Value = -21.74339999999999
$ws.Range("B78").Value = 5.9207
$ws.Range("B80").Value = 9.416699999999997
$ws.Range("A82").Value = -21.99290000000001
$ws.Range("C82").Value = -11.3465
$ws.Range("A83").Value = -21.57019999999998
$ws.Range("C90").Value = -10.16070000000001
$ws.Range("B91").Value = 7.591900000000003
$ws.Range("D91").Value = -7.489399999999998
$ws.Range("D92").Value = -6.350800000000001
$ws.Range("A93").Value = -21.55610000000001
$ws.Range("C93").Value = -11.0084
$ws.Range("D95").Value = -7.462600000000004
$ws.Range("D96").Value = -7.794299999999994
$ws.Range("A97").Value = -21.5624
$ws.Range("B97").Value = 4.794099999999991
$ws.Range("B99").Value = 6.2834
$ws.Range("C102").Value = -11.9512
$ws.Range("D103").Value = -8.598099999999997
$ws.Range("B104").Value = 9.985200000000001
$ws.Range("C105").Value = -12.40170000000001
$ws.Range("D105").Value = -7.542199999999999
